# Fix Week 12 date: change "12/20" to "11/20" in both the Week label (col B)
# and the Date column (col C), for every student's Week 12 row.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$rows = 5..233 | Where-Object { ($_ - 5) % 12 -eq 0 }

foreach ($r in $rows) {
    $ws.Range("B$r").Value = "Week 12 (11/20)"
    $ws.Range("C$r").Value = "11/20"
}
